$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S1").Value = "ADJ1"
$ws.Range("T1").Value = "ADJ2"
$ws.Range("U1").Value = "ADJ3"

$ws.Range("S1:U1").Select()
